# Updates the Price (D) and Volume(1h) (E) columns of the cryptos
# worksheet to the latest scraped snapshot.
#
# For D-column values that would otherwise be auto-parsed by Excel as a
# number (single "." decimal separator, e.g. "236.89" or "0.0000303"),
# the cell's NumberFormat is forced to "@" (Text) immediately before the
# assignment so the literal text is preserved exactly (matching the
# source, which keeps these as plain text cells) instead of being
# silently coerced into a numeric value/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.111.80'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '3.694.71'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '236.89'
$ws.Range('E5').Value = '  -2.08%  '
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '658.63'
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.424'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.999'
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '3.693.98'
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '44.20'
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000303'
$ws.Range('E14').Value = '  +11.47%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.77'
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').Value = '4.382.07'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '97.043.46'
$ws.Range('E17').Value = '  +0.89%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '9.14'
$ws.Range('E18').Value = '  +2.49%  '
$ws.Range('D19').Value = '3.746.72'
$ws.Range('E19').Value = '  +2.44%  '
$ws.Range('E20').Value = '  +2.13%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '18.77'
$ws.Range('E21').Value = '  +2.79%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.509'
$ws.Range('E22').Value = '  -4.84%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '520.27'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('E25').Value = '  +2.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.94'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('E27').Value = '  +22.86%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '101.47'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E29').Value = '  +3.26%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '12.53'
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('E34').Value = '  +2.28%  '
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '32.20'
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '647.26'
$ws.Range('E37').Value = '  +3.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.592'
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('E39').Value = '  +1.47%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.504'
$ws.Range('E41').Value = '  +13.94%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.84'
$ws.Range('E42').Value = '  +9.05%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.06'
$ws.Range('E43').Value = '  +5.37%  '
$ws.Range('E44').Value = '  +1.36%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '40.46'
$ws.Range('E45').Value = '  -10.53%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.961'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0467'
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.28'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.71'
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('E51').Value = '  -1.24%  '
